# Update "F" column (想去人数 / want-to-go count) figures across sheets
# to reflect newly generated output (gh-pages regeneration at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 712
    3  = 60
    4  = 2017
    5  = 5911
    6  = 1668
    7  = 182
    8  = 3352
    11 = 1404
    12 = 4660
    13 = 1759
    14 = 19
    15 = 59
    16 = 62
    17 = 208
    19 = 1047
    20 = 320
    22 = 29
    23 = 93
    27 = 1140
    28 = 427
    29 = 109
    30 = 224
    31 = 462
    34 = 1800
    35 = 2286
    36 = 1072
    38 = 5
    40 = 8
    41 = 651
    42 = 420
    43 = 50
    45 = 43
    46 = 460
    47 = 445
    49 = 151
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# --- Sheet "本地生活" (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 805

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 712
    3  = 60
    4  = 2017
    5  = 5911
    6  = 1668
    7  = 182
    9  = 3352
    11 = 1404
    12 = 4660
    13 = 1759
    14 = 19
    16 = 59
    20 = 62
    21 = 208
    24 = 1047
    25 = 320
    27 = 29
    28 = 93
    32 = 1140
    33 = 109
    34 = 224
    37 = 1800
    38 = 2286
    47 = 460
    48 = 445
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
